# Apply "Updated 19 test cases with loop loading details method" change
# to the "Add Devices Loop A" worksheet (first sheet in the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New header / data cells in column E (Volt Drop (V) / Volt Drop (worst case)) ---
# These two new cells take on the same shaded/bordered look as the existing
# green/gray data cells (e.g. A6) but left-aligned with wrap text.

$ws.Range("E3").Value = "Volt Drop (V)"
$ws.Range("E4").Value = "Volt Drop (worst case)"

$ws.Range("A6").Copy() | Out-Null
$ws.Range("E3:E4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("E3:E4").HorizontalAlignment = -4131       # xlLeft
$ws.Range("E3:E4").WrapText = $true

# --- New header cells H5 / I5 (Loading Details Name headers) ---
# Same bold / bordered header look as the existing G5 header cell.

$ws.Range("H5").Value = "Volt drop Loading Details Name"
$ws.Range("I5").Value = "Volt drop worst case Loading Details Name"

$ws.Range("G5").Copy() | Out-Null
$ws.Range("H5:I5").PasteSpecial(-4122) | Out-Null    # xlPasteFormats
$excel.CutCopyMode = 0

# --- New data cells H6:I8 (Loading Details Name values, one per device row) ---
# Same shaded/bordered look as the other data cells in that row (e.g. A6).

$ws.Range("H6").Value = "Volt Drop (V)"
$ws.Range("I6").Value = "Volt Drop (worst case)"
$ws.Range("H7").Value = "Volt Drop (V)"
$ws.Range("I7").Value = "Volt Drop (worst case)"
$ws.Range("H8").Value = "Volt Drop (V)"
$ws.Range("I8").Value = "Volt Drop (worst case)"

$ws.Range("A6").Copy() | Out-Null
$ws.Range("H6:I8").PasteSpecial(-4122) | Out-Null    # xlPasteFormats
$excel.CutCopyMode = 0

# --- Column widths: column E is new, column H got wider to fit new text ---
$ws.Columns.Item(5).ColumnWidth = 26.33203125
$ws.Columns.Item(8).ColumnWidth = 27.88671875

# --- Update the active selection to match the edited range ---
$ws.Range("E3:E4").Select() | Out-Null
